# Update classification result metrics (Cross Entropy Loss / Success % columns,
# plus a couple of refreshed Prediction/Error values) to reflect the newly
# regenerated toy NCDE/NODE datasets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Batch size 100 block (rows 2-11) ---
$ws.Range("D2").Value = 0.0007754405458628474
$ws.Range("E2").Value = 0.0007754405458628474

$ws.Range("D3").Value = 0.8744353818827404
$ws.Range("E3").Value = 0.8744353818827404

$ws.Range("D4").Value = 0.001072861189595828
$ws.Range("E4").Value = 0.001072861189595828

$ws.Range("D5").Value = 0.0000000000005695518310458689
$ws.Range("E5").Value = 0.0000000000005695518310458689

$ws.Range("D6").Value = 0.7221888631987826
$ws.Range("E6").Value = 0.7221888631987826

$ws.Range("D7").Value = 0.748266836408498
$ws.Range("E7").Value = 0.251733163591502

$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.0000001218166366293734
$ws.Range("E8").Value = 0.9999998781833633

$ws.Range("D9").Value = 0.9418445522083554
$ws.Range("E9").Value = 0.05815544779164461

$ws.Range("D10").Value = 0.5491940909700388
$ws.Range("E10").Value = 0.4508059090299612

$ws.Range("D11").Value = 0.9999999916341885
$ws.Range("E11").Value = 0.000000008365811487998087
$ws.Range("F11").Value = 2.022756099700928
$ws.Range("G11").Value = 0.7

# --- Batch size 200 block (rows 12-21) ---
$ws.Range("D12").Value = 0.00000007222562763807043
$ws.Range("E12").Value = 0.00000007222562763807043

$ws.Range("D13").Value = 0.8746887032394338
$ws.Range("E13").Value = 0.8746887032394338

$ws.Range("D14").Value = 0.000000146192447132351
$ws.Range("E14").Value = 0.000000146192447132351

$ws.Range("D15").Value = 0.000000000002667342146301273
$ws.Range("E15").Value = 0.000000000002667342146301273

$ws.Range("D16").Value = 0.9597346355012241
$ws.Range("E16").Value = 0.9597346355012241

$ws.Range("D17").Value = 0.9555242470885883
$ws.Range("E17").Value = 0.0444757529114117

$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.000001167333019782722
$ws.Range("E18").Value = 0.9999988326669802

$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.1285291646841141
$ws.Range("E19").Value = 0.871470835315886

$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.000004450950293059468
$ws.Range("E20").Value = 0.9999955490497069

$ws.Range("D21").Value = 0.999999999999974
$ws.Range("E21").Value = 0.00000000000002597921877622866
$ws.Range("F21").Value = 3.336949586868286
$ws.Range("G21").Value = 0.5
